# Word COM-interop script implementing the "04. Glossário.docx" edits:
#   1. "Automova" -> "Automotiva"               (Pintura Automo[va->tiva] row)
#   2. "azer um pintura" -> "azer uma pintura"   (Fazer um[+a] pintura no automóvel)
#
# (The source diff also re-wraps an already-correct "politriz," span with
#  <w:proofErr> spell-check markers and relocates the hidden "_GoBack" last-
#  edit-position bookmark; those are Word's own internal, invisible
#  bookkeeping artifacts written by the live spell-checker / caret tracker,
#  not user-addressable content, and they do not alter the document text.)

$d = $word.ActiveDocument

# --- 1. "Pintura Automova" -> "Pintura Automotiva" ---------------------
$found1 = $d.Content.Find.Execute(
    "Automova", $true, $false, $false, $false, $false,
    $true, 1, $false, "Automotiva", 2)
if (-not $found1) {
    throw "Could not find 'Automova' to replace with 'Automotiva'"
}

# --- 2. "Fazer um pintura no automóvel" -> "Fazer uma pintura no automóvel" ---
$found2 = $d.Content.Find.Execute(
    "azer um pintura no automóvel", $true, $false, $false, $false, $false,
    $true, 1, $false, "azer uma pintura no automóvel", 2)
if (-not $found2) {
    throw "Could not find 'azer um pintura no automóvel' to replace"
}
